$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7993
$ws.Range("I64").Value = 7989
$ws.Range("K64").Value = 7989
$ws.Range("M64").Value = -7741
$ws.Range("H67").Value = 7993
$ws.Range("I67").Value = 7989
$ws.Range("K67").Value = 7989
$ws.Range("M67").Value = -7131
$ws.Range("H96").Value = 1246.8823
$ws.Range("I96").Value = 1342.7142
$ws.Range("K96").Value = 4028.1426
$ws.Range("M96").Value = -2655.1426
$ws.Range("H113").Value = 3116.6667
$ws.Range("I113").Value = 3175
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 3175
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 79
$ws.Range("N113").Value = -9508
$ws.Range("H138").Value = 4088.4375
$ws.Range("I138").Value = 3348.25
$ws.Range("J138").Value = 4335.1665
$ws.Range("K138").Value = 10044.75
$ws.Range("L138").Value = 13005.4995
$ws.Range("M138").Value = -4904.75
$ws.Range("N138").Value = -23285.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3302.6667
$ws.Range("I32").Value = 3257.6316
$ws.Range("K32").Value = 3257.6316
$ws.Range("M32").Value = -2970.6316
$ws.Range("H37").Value = 12333
$ws.Range("I37").Value = 7499.5
$ws.Range("J37").Value = 22000
$ws.Range("K37").Value = 7499.5
$ws.Range("L37").Value = 22000
$ws.Range("M37").Value = -7226.5
$ws.Range("N37").Value = -22546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2363.3333
$ws.Range("I80").Value = 200
$ws.Range("J80").Value = 2796
$ws.Range("K80").Value = 200
$ws.Range("L80").Value = 2796
$ws.Range("M80").Value = 798
$ws.Range("N80").Value = -4792
$ws.Range("H83").Value = 2363.3333
$ws.Range("I83").Value = 200
$ws.Range("J83").Value = 2796
$ws.Range("K83").Value = 1000
$ws.Range("L83").Value = 13980
$ws.Range("M83").Value = 3992
$ws.Range("N83").Value = -23964

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H62").Value = 4003.8
$ws.Range("I62").Value = 3753.25
$ws.Range("J62").Value = 5006
$ws.Range("K62").Value = 3753.25
$ws.Range("L62").Value = 5006
$ws.Range("M62").Value = -3129.25
$ws.Range("N62").Value = -6254
$ws.Range("H65").Value = 4003.8
$ws.Range("I65").Value = 3753.25
$ws.Range("J65").Value = 5006
$ws.Range("K65").Value = 18766.25
$ws.Range("L65").Value = 25030
$ws.Range("M65").Value = -15646.25
$ws.Range("N65").Value = -31270
$ws.Range("H105").Value = 2853
$ws.Range("I105").Value = 2880
$ws.Range("J105").Value = 2691
$ws.Range("K105").Value = 2880
$ws.Range("L105").Value = 2691
$ws.Range("M105").Value = -1133
$ws.Range("N105").Value = -6185
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 99.5
$ws.Range("H8").Value = 999
$ws.Range("I8").Value = 999
$ws.Range("K8").Value = 2997
$ws.Range("M8").Value = -2858
$ws.Range("H15").Value = 449
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 449
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 1347
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -1627
$ws.Range("H17").Value = 212.5
$ws.Range("I17").Value = 212.5
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 637.5
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -468.5
$ws.Range("N17").ClearContents()
$ws.Range("H18").Value = 3058.2856
$ws.Range("I18").Value = 1061.6
$ws.Range("J18").Value = 8050
$ws.Range("K18").Value = 3184.8
$ws.Range("L18").Value = 24150
$ws.Range("M18").Value = -3015.8
$ws.Range("N18").Value = -24488
$ws.Range("H34").Value = 3750.625
$ws.Range("J34").Value = 3997.5833
$ws.Range("L34").Value = 11992.7499
$ws.Range("N34").Value = -12160.7499
$ws.Range("H39").Value = 16000
$ws.Range("J39").Value = 16000
$ws.Range("L39").Value = 48000
$ws.Range("N39").Value = -48588
$ws.Range("H55").Value = 11307.714
$ws.Range("J55").Value = 15511
$ws.Range("L55").Value = 46533
$ws.Range("N55").Value = -46887
$ws.Range("H121").Value = 1249.5
$ws.Range("I121").Value = 699
$ws.Range("K121").Value = 2097
$ws.Range("M121").Value = -787
$ws.Range("H139").Value = 2209.5
$ws.Range("I139").Value = 2209.5
$ws.Range("K139").Value = 6628.5
$ws.Range("M139").Value = -1488.5
$ws.Range("H140").Value = 5862.8
$ws.Range("I140").Value = 5817.067
$ws.Range("K140").Value = 17451.201
$ws.Range("M140").Value = -12271.201
$ws.Range("H141").Value = 963.3333
$ws.Range("I141").Value = 963.3333
$ws.Range("K141").Value = 2889.9999
$ws.Range("M141").Value = 2290.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5193.375
$ws.Range("J22").Value = 7339.9
$ws.Range("L22").Value = 7339.9
$ws.Range("N22").Value = -7929.9
$ws.Range("H27").Value = 5193.375
$ws.Range("J27").Value = 7339.9
$ws.Range("L27").Value = 7339.9
$ws.Range("N27").Value = -7553.9
$ws.Range("H40").Value = 3066.8333
$ws.Range("I40").Value = 3080.2
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 3080.2
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2944.2
$ws.Range("N40").Value = -3272

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 30000
$ws.Range("J54").Value = 30000
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -31040
$ws.Range("H81").Value = 12004.5
$ws.Range("I81").Value = 12004.5
$ws.Range("K81").Value = 24009
$ws.Range("M81").Value = -22948
$ws.Range("H84").Value = 12004.5
$ws.Range("I84").Value = 12004.5
$ws.Range("K84").Value = 120045
$ws.Range("M84").Value = -114741
